# Assignment 1 report edit: append a new explanatory sentence to the end
# of the paragraph describing the "fget" function in the C++ Side of the
# Bonus Deliverable (Dialogue System) section.
#
# Before: "...skipping lines until it reaches the desired line number as
#          specified in the integer variable."
# After:  "...skipping lines until it reaches the desired line number as
#          specified in the integer variable. When it reaches the desired
#          line, it stores the information in that and each subsequent
#          line into the stringstream variable, until it reaches a line
#          that says >BREAK<, which forces the while loop to end early."

$d = $word.ActiveDocument

$oldText = "as specified in the integer variable."
$newText = "as specified in the integer variable. When it reaches the desired line, it stores the information in that and each subsequent line into the stringstream variable, until it reaches a line that says >BREAK<, which forces the while loop to end early."

$found = $d.Content.Find.Execute(
    $oldText, $true, $true, $false, $false, $false,
    $true, 1, $false, $newText, 2)

if (-not $found) {
    Write-Host "WARNING: target sentence not found; no replacement made."
} else {
    Write-Host "Inserted continuation sentence after the fget description."
}
